$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CREAR TABLAS")

# Fill in the "SCRIPT" column (D) for rows whose script was just written, and
# bump their progress column (E) from 0 to 10 — matching the "Scripts
# Catalogo y tablas" commit.
#
# Order matters: new shared strings are appended in first-use order, so we
# touch the rows in the same sequence the diff implies (9, 33, 41, 45, 67,
# then 53) to reproduce the same sharedStrings.xml layout.

$ws.Range("D9").Value = "T2-Tabla_Email"
$ws.Range("E9").Value = 10

$ws.Range("D33").Value = "C26-Catalogo_Ciudad"
$ws.Range("E33").Value = 10

$ws.Range("D41").Value = "C34-Catalogo_Contextura"
$ws.Range("E41").Value = 10

$ws.Range("D45").Value = "C38-Catalogo_Cojos"
$ws.Range("E45").Value = 10

$ws.Range("D67").Value = "C60-Catalogo_Escaolaridad"
$ws.Range("E67").Value = 10

$ws.Range("D53").Value = "C46-Catalogo_Ocupacion"
$ws.Range("E53").Value = 10

# Move the view / active selection the way the author left it.
$ws.Range("D54").Select()
